$d = $word.ActiveDocument

# Use Find/Replace to split the name paragraph into two paragraphs:
# the existing "Dheeraj Chand" line, followed by a new centered
# paragraph containing the contact information. Replacing this way
# (rather than Range.InsertParagraphAfter) avoids carrying the bold /
# large-font run formatting of the name onto the new contact line.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Dheeraj Chand", $true, $false, $false, $false, $false, $true, 1, $false,
              "Dheeraj Chand^p202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX",
              2)
